{"js": "// 1) Make the \"Bibliografie\" title bold and bump its font size from 12pt to 16pt.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst title = paragraphs.items[0];\ntitle.font.bold = true;\ntitle.font.size = 16;\n\n// 2) Remove the stray \"Reference :\" text that prefixed the first bibliography entry.\nconst hits = context.document.body.search(\"Reference :\", { matchCase: false });\nhits.load(\"items\");\nawait context.sync();\n\nif (hits.items.length > 0) {\n  hits.items[0].delete();\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) Make the \"Bibliografie\" title bold and bump its font size from 12pt to 16pt\n#    (regular + complex-script formatting, matching the diacritics/complex-script aware\n#    edit described in the commit message).\n$title = $d.Paragraphs(1).Range\n$title.Font.Bold = 1\n$title.Font.BoldBi = 1\n$title.Font.Size = 16\n$title.Font.SizeBi = 16\n\n# 2) Remove the stray \"Reference :\" text that prefixed the first bibliography entry.\n$find = $d.Content\n$find.Find.ClearFormatting()\n$find.Find.Text = \"Reference :\"\n$find.Find.MatchCase = $false\n$find.Find.Forward = $true\n$find.Find.Wrap = 1\nif ($find.Find.Execute()) {\n    $find.Delete()\n}\n"}
